$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")
$rng = $ws.Range("A8")
$rng.Value = "discord"
$rng.Font.ThemeFont = 2
$rng.Font.Bold = $true
Write-Host "done"
